$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2:E51 hold price/volume figures as plain text (e.g. "59.468.37", "  +2.72%  ").
# Force the range to Text format before writing so Excel does not reinterpret
# strings like "5.38" or "35.20" as numbers, then restore the original style.
$savedStyle = $ws.Range("D2:E51").Style
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "59.468.37"
$ws.Range("E2").Value = "  +2.72%  "

$ws.Range("D3").Value = "2.984.13"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "566.67"
$ws.Range("E5").Value = "  +2.27%  "

$ws.Range("D6").Value = "138.42"
$ws.Range("E6").Value = "  +3.83%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +1.36%  "

$ws.Range("D9").Value = "2.973.31"
$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("E10").Value = "  +3.32%  "

$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  +11.43%  "

$ws.Range("D12").Value = "0.451"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("E13").Value = "  +3.65%  "

$ws.Range("D14").Value = "33.73"
$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").Value = "3.477.09"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("D17").Value = "7.03"
$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").Value = "2.979.97"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("D19").Value = "59.468.00"
$ws.Range("E19").Value = "  +2.80%  "

$ws.Range("D20").Value = "436.79"
$ws.Range("E20").Value = "  +4.63%  "

$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("E22").Value = "  +3.04%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  -1.34%  "

$ws.Range("D25").Value = "79.98"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +9.91%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  +2.07%  "

$ws.Range("D30").Value = "7.74"
$ws.Range("E30").Value = "  +2.87%  "

$ws.Range("D31").Value = "25.74"
$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("D32").Value = "6.21"
$ws.Range("E32").Value = "  +3.90%  "

$ws.Range("E33").Value = "  +8.77%  "

$ws.Range("D34").Value = "0.0₃0774"
$ws.Range("E34").Value = "  +10.70%  "

$ws.Range("D35").Value = "5.90"
$ws.Range("E35").Value = "  +3.60%  "

$ws.Range("E36").Value = "  +3.64%  "

$ws.Range("D37").Value = "2.08"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").Value = "48.61"
$ws.Range("E38").Value = "  +0.53%  "

$ws.Range("E39").Value = "  -3.72%  "

$ws.Range("D40").Value = "2.78"
$ws.Range("E40").Value = "  +2.56%  "

$ws.Range("D41").Value = "400.96"
$ws.Range("E41").Value = "  +4.90%  "

$ws.Range("D42").Value = "0.0351"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").Value = "2.734.37"
$ws.Range("E43").Value = "  +0.98%  "

$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("E45").Value = "  +5.34%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "35.20"
$ws.Range("E46").Value = "  +19.73%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").Value = "121.93"
$ws.Range("E48").Value = "  -1.91%  "

$ws.Range("E49").Value = "  +1.58%  "

$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("D51").Value = "23.29"
$ws.Range("E51").Value = "  +1.46%  "

$ws.Range("D2:E51").Style = $savedStyle
